# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gains a new (blank) column inserted just
# before the existing "Late" column (column N), pushing "Late",
# "heading"/Disbursement and "Outstanding" one column to the right
# (N->O, O->P, P->Q). The sheet becomes the active tab, with R10 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column immediately to the left (M) so the
# newly inserted column inherits it, matching Excel's native insert
# behaviour.
$leftWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before column N ("Late"), shifting the
# "Late", heading and "Outstanding" columns one place to the right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Make "Repayment schedule" the active sheet/tab and select cell R10.
$ws.Activate() | Out-Null
$ws.Range("R10").Select() | Out-Null
